# Updated cryptos list (GitHub Actions refresh of prices / 1h volume%).
# Two rows (20/21, 33/35/36, 45/46) also had their Coin/Link/Price/Volume
# swapped with a neighboring row as the ranking shifted.
# Numeric-looking Price strings are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the source inlineStr cells) instead
# of auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.800.28"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.565.45"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'311.61"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'98.29"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "'35.61"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "'7.42"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").Value = "2.967.23"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").Value = "'15.86"
$ws.Range("E15").Value = "  +4.58%  "
$ws.Range("D16").Value = "2.525.95"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "'0.840"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "42.849.04"
$ws.Range("D19").Value = "'6.73"
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0958"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "'12.39"
$ws.Range("E21").Value = "  -4.17%  "
$ws.Range("D22").Value = "'69.77"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "'248.38"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").Value = "'2.93"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'2.05"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "'27.00"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "'39.63"
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("D30").Value = "'10.21"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "'159.16"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").Value = "'5.77"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.68"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").Value = "'2.10"
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'3.31"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0795"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").Value = "'18.53"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").Value = "'2.57"
$ws.Range("E38").Value = "  +10.61%  "
$ws.Range("D39").Value = "'0.111"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Value = "'22.84"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").Value = "'4.11"
$ws.Range("E42").Value = "  +7.14%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "'0.0301"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.996.02"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'3.20"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("D47").Value = "'8.99"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").Value = "2.818.48"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "'0.195"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("D50").Value = "'81.55"
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("D51").Value = "'74.04"
$ws.Range("E51").Value = "  -0.85%  "
